$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column I entirely (was epoch400 header + its value)
$ws.Range("I1:I2").Delete()

# Update the values in row 2 with newly trained model results
$ws.Range("B2").Value = 48.32175903850131
$ws.Range("C2").Value = 38.83101807700263
$ws.Range("D2").Value = 34.02777777777778
$ws.Range("E2").Value = 27.60416666666666
$ws.Range("F2").Value = 25.92592570516798
$ws.Range("G2").Value = 25.23148126072354
$ws.Range("H2").Value = 25.23148126072354
